$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Markers")

# Column E (rows 3-9) becomes an exact copy (value + style) of column D
$ws.Range("D3").Copy($ws.Range("E3"))
$ws.Range("D4").Copy($ws.Range("E4"))
$ws.Range("D5").Copy($ws.Range("E5"))
$ws.Range("D6").Copy($ws.Range("E6"))
$ws.Range("D7").Copy($ws.Range("E7"))
$ws.Range("D8").Copy($ws.Range("E8"))
$ws.Range("D9").Copy($ws.Range("E9"))

# Column E (rows 10-17) picks up column D's style, but stays empty
$ws.Range("D10").Copy($ws.Range("E10"))
$ws.Range("E10").ClearContents()
$ws.Range("D11").Copy($ws.Range("E11"))
$ws.Range("E11").ClearContents()
$ws.Range("D12").Copy($ws.Range("E12"))
$ws.Range("E12").ClearContents()
$ws.Range("D13").Copy($ws.Range("E13"))
$ws.Range("E13").ClearContents()
$ws.Range("D14").Copy($ws.Range("E14"))
$ws.Range("E14").ClearContents()
$ws.Range("D15").Copy($ws.Range("E15"))
$ws.Range("E15").ClearContents()
$ws.Range("D16").Copy($ws.Range("E16"))
$ws.Range("E16").ClearContents()
$ws.Range("D17").Copy($ws.Range("E17"))
$ws.Range("E17").ClearContents()

# Update the saved selection to match the committed state
$ws.Range("I17").Select()
